$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column Q ("EDAD"): age in whole years computed from birth date (C) to today
$ws.Range("Q4").Formula = "=INT(YEARFRAC(C4,TODAY()))"
$ws.Range("Q5:Q32").Formula = "=INT(YEARFRAC(C5,TODAY()))"

# Column R ("VALIMENTACIÓN"): meal allowance based on net salary (O) and hourly-rate flag (F)
$ws.Range("R4").Formula = "=IF(O4<5000,300,IF(AND(O4>=5000,O4<=10000),IF(ISNUMBER(F4),200,100),0))"
$ws.Range("R5:R32").Formula = "=IF(O5<5000,300,IF(AND(O5>=5000,O5<=10000),IF(ISNUMBER(F5),200,100),0))"

# Move the active selection to R8
$ws.Range("R8").Select() | Out-Null
